$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store their figures as literal
# text (e.g. "312.00", "0.46%"), not numbers. Excel's COM Value setter
# auto-converts plain numeric-looking / percent-looking strings, so each
# target cell is switched to Text format ("@") before its new value is
# written, to preserve the exact text representation from the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "312.00"
$ws.Range("E2").Value = "0.46%"
$ws.Range("D3").Value = "38.48"
$ws.Range("E3").Value = "-2.46%"
$ws.Range("D4").Value = "5.133"
$ws.Range("E4").Value = "0.22%"
$ws.Range("D5").Value = "0.08101"
$ws.Range("E5").Value = "-0.09%"
$ws.Range("D6").Value = "4.469"
$ws.Range("E6").Value = "5.49%"
$ws.Range("D7").Value = "1.954"
$ws.Range("E7").Value = "-3.27%"
$ws.Range("D8").Value = "8.306"
$ws.Range("E8").Value = "1.88%"
$ws.Range("D9").Value = "0.9390"
$ws.Range("E9").Value = "1.01%"
$ws.Range("D10").Value = "0.1326"
$ws.Range("E10").Value = "-7.26%"
$ws.Range("D11").Value = "0.1937"
$ws.Range("E11").Value = "0.29%"
$ws.Range("D12").Value = "0.09035"
$ws.Range("E12").Value = "-0.56%"
$ws.Range("D13").Value = "0.03485"
$ws.Range("E13").Value = "-0.43%"
$ws.Range("D14").Value = "0.09673"
$ws.Range("E14").Value = "-1.48%"
$ws.Range("D15").Value = "0.001407"
$ws.Range("E15").Value = "0.35%"
$ws.Range("D16").Value = "0.005888"
$ws.Range("E16").Value = "-0.04%"
$ws.Range("D17").Value = "3.566"
$ws.Range("E17").Value = "-6.68%"
$ws.Range("D18").Value = "3.191"
$ws.Range("E18").Value = "-4.17%"
$ws.Range("D19").Value = "0.3464"
$ws.Range("E19").Value = "0.32%"
$ws.Range("E20").Value = "-3.84%"
$ws.Range("D21").Value = "5.012"
$ws.Range("E21").Value = "6.02%"
$ws.Range("D22").Value = "0.2498"
$ws.Range("E22").Value = "2.91%"
$ws.Range("D23").Value = "0.04368"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").Value = "0.38%"
$ws.Range("D25").Value = "0.004728"
$ws.Range("E25").Value = "-1.44%"
$ws.Range("D26").Value = "0.0003898"
$ws.Range("E26").Value = "199.36%"
$ws.Range("D39").Value = "0.02203"
$ws.Range("E39").Value = "3.46%"
$ws.Range("D40").Value = "0.05228"
$ws.Range("E40").Value = "2.54%"
$ws.Range("D41").Value = "0.007600"
$ws.Range("E41").Value = "2.03%"
$ws.Range("D42").Value = "0.01035"
$ws.Range("E42").Value = "5.32%"
$ws.Range("D43").Value = "0.1389"
$ws.Range("E43").Value = "2.04%"
$ws.Range("D44").Value = "0.002033"
$ws.Range("E44").Value = "-4.72%"
$ws.Range("D45").Value = "0.009104"
$ws.Range("E45").Value = "5.61%"
$ws.Range("D46").Value = "0.00006614"
$ws.Range("E46").Value = "3.48%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D48").Value = "0.003014"
$ws.Range("E48").Value = "16.96%"
$ws.Range("E49").Value = "68.85%"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.04%"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "0.04%"
